$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "65.06"
$ws.Range("B13").NumberFormat = "General"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "2.86"
$ws.Range("C13").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "67.92"
$ws.Range("D13").NumberFormat = "General"

# Employment (% of total): row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "30.22"
$ws.Range("B14").NumberFormat = "General"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "32.62"
$ws.Range("C14").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "62.84"
$ws.Range("D14").NumberFormat = "General"

# Enterprises (% of total): row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "93.79"
$ws.Range("B16").NumberFormat = "General"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "6.08"
$ws.Range("C16").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.88"
$ws.Range("D16").NumberFormat = "General"
